{"js": "// Find the \"HSTS : forcez le passage \u00e0 HTTPS et www\" list item and insert a\n// new list item \"Cache config\" immediately after it (same list / level).\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (t.indexOf(\"HSTS\") !== -1 && t.indexOf(\"HTTPS\") !== -1) {\n    target = paras.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the \"HSTS\" paragraph');\n}\n\n// insertParagraph inherits the source paragraph's style and list\n// membership (numId/ilvl), which is exactly what we want here.\ntarget.insertParagraph(\"Cache config\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Find the \"HSTS : forcez le passage \u00e0 HTTPS et www\" list item and insert a\n# new list item \"Cache config\" immediately after it (same list / level).\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*HSTS*HTTPS*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the HSTS paragraph\"\n}\n\n# Collapse to the end of the paragraph (past the paragraph mark) and insert\n# a new paragraph there; the new paragraph inherits the style/list (numId,\n# ilvl) of the paragraph it follows.\n$r = $target.Range\n$r.Collapse(0)  # wdCollapseEnd\n$r.InsertParagraphAfter()\n\n$newRange = $target.Next().Range\n$newRange.Text = \"Cache config\"\n"}
